$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K42").Value = "test"
